# Updates the cryptos price (D) / volume-change (E) columns and swaps the
# Hedera / LidoDAOToken rows (35 <-> 36), matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "54.48", "1.00", "0.0800")
# get a leading apostrophe so Excel keeps storing them as text (preserving
# exact digits/trailing zeros) instead of silently converting to a numeric value.

$ws.Range("D2").Value = "41.886.82"
$ws.Range("E2").Value = "  +5.19%  "
$ws.Range("D3").Value = "2.271.80"
$ws.Range("D5").Value = "'303.15"
$ws.Range("E5").Value = "  +4.16%  "
$ws.Range("D6").Value = "'93.15"
$ws.Range("E6").Value = "  +7.76%  "
$ws.Range("E7").Value = "  +3.16%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +5.62%  "
$ws.Range("D10").Value = "'54.48"
$ws.Range("E10").Value = "  +8.79%  "
$ws.Range("D11").Value = "'32.43"
$ws.Range("E11").Value = "  +8.50%  "
$ws.Range("D12").Value = "'0.0800"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").Value = "'6.68"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").Value = "2.623.83"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("D16").Value = "'14.21"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").Value = "2.261.35"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "'0.756"
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("D19").Value = "41.742.36"
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("D20").Value = "'12.51"
$ws.Range("E20").Value = "  +12.25%  "
$ws.Range("D21").Value = "0.0₃0912"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("E22").Value = "  +4.30%  "
$ws.Range("D23").Value = "'67.23"
$ws.Range("E23").Value = "  +3.46%  "
$ws.Range("D24").Value = "'241.01"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "'2.59"
$ws.Range("E25").Value = "  +6.62%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("E29").Value = "  +7.05%  "
$ws.Range("E30").Value = "  +6.84%  "
$ws.Range("D31").Value = "'34.29"
$ws.Range("E31").Value = "  +10.23%  "
$ws.Range("D32").Value = "'158.23"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +6.77%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0739"
$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.07"
$ws.Range("E36").Value = "  +9.90%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").Value = "'16.67"
$ws.Range("E38").Value = "  +10.88%  "
$ws.Range("E39").Value = "  +7.54%  "
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("E41").Value = "  +8.13%  "
$ws.Range("D42").Value = "'3.99"
$ws.Range("E42").Value = "  +7.47%  "
$ws.Range("D43").Value = "'20.37"
$ws.Range("E43").Value = "  +19.51%  "
$ws.Range("D44").Value = "2.063.21"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("E46").Value = "  +12.80%  "
$ws.Range("D47").Value = "'10.07"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").Value = "'1.99"
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D49").Value = "2.494.22"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("E50").Value = "  +4.82%  "
$ws.Range("D51").Value = "'1.14"
$ws.Range("E51").Value = "  +4.99%  "
